$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New test cases to append after row 57, matching the plain (unshaded)
# row formatting used by rows like 43-47.
$ws.Range("A47:D47").Copy() | Out-Null
$ws.Range("A58:D60").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item(58, 1).Value = "EGN36"
$ws.Cells.Item(58, 2).Value = "cmd 81-0C/ 81-0D"
$ws.Cells.Item(58, 3).Value = ""
$ws.Cells.Item(58, 4).Value = "DONE"

$ws.Cells.Item(59, 1).Value = "EGN37"
$ws.Cells.Item(59, 2).Value = "CS-5449, NSRED reader, cmd C7-3A w/ 01 (internal cmd)"
$ws.Cells.Item(59, 3).Value = ""
$ws.Cells.Item(59, 4).Value = "DONE"

$ws.Cells.Item(60, 1).Value = "EGN38"
$ws.Cells.Item(60, 2).Value = "Multiple DEK test, ex. KeySlot 00 = TDES KEY, KeySlot 03 = AES KEY"
$ws.Cells.Item(60, 3).Value = ""
$ws.Cells.Item(60, 4).Value = "DONE"

$ws.Range("B62").Select() | Out-Null
